$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.920821042693349
$ws.Range("D2").Value = 9.049374931987233
$ws.Range("E2").Value = 14.87132837618671
$ws.Range("F2").Value = 23.57344825814363
$ws.Range("G2").Value = 22.86492872204558
$ws.Range("H2").Value = 12.71982417213049
$ws.Range("J2").Value = 10.12527413481625
$ws.Range("M2").Value = 59.02342950016224
$ws.Range("O2").Value = 18.58332014747845
$ws.Range("C3").Value = 4.749564180890503
$ws.Range("D3").Value = 9.088453672064107
$ws.Range("E3").Value = 14.59835990872851
$ws.Range("F3").Value = 24.06004497121543
$ws.Range("G3").Value = 23.17267786785296
$ws.Range("H3").Value = 12.85522549006706
$ws.Range("J3").Value = 10.01032181506775
$ws.Range("M3").Value = 55.55207178020056
$ws.Range("O3").Value = 18.82296337929108
$ws.Range("C4").Value = 4.642847137742884
$ws.Range("D4").Value = 9.115954648386445
$ws.Range("E4").Value = 14.4352822732601
$ws.Range("F4").Value = 24.37335300164089
$ws.Range("G4").Value = 23.38559645693615
$ws.Range("H4").Value = 12.94303752811742
$ws.Range("J4").Value = 9.943375437464914
$ws.Range("M4").Value = 53.29776787465694
$ws.Range("O4").Value = 18.98039990439404
$ws.Range("C5").Value = 4.599042422035566
$ws.Range("D5").Value = 9.128032462041402
$ws.Range("E5").Value = 14.3700502232123
$ws.Range("F5").Value = 24.50464912464627
$ws.Range("G5").Value = 23.4781854778397
$ws.Range("H5").Value = 12.97998898397182
$ws.Range("J5").Value = 9.917030498814444
$ws.Range("M5").Value = 52.34840256380174
$ws.Range("O5").Value = 19.0471028815053
$ws.Range("C6").Value = 4.591751920823855
$ws.Range("D6").Value = 9.130090266308651
$ws.Range("E6").Value = 14.35929484925585
$ws.Range("F6").Value = 24.52666847233554
$ws.Range("G6").Value = 23.49390570592758
$ws.Range("H6").Value = 12.98619501584465
$ws.Range("J6").Value = 9.91271310744359
$ws.Range("M6").Value = 52.18891252212197
$ws.Range("O6").Value = 19.05833145014889
$ws.Range("C7").Value = 4.642257548175663
$ws.Range("D7").Value = 9.116114021626833
$ws.Range("E7").Value = 14.43439746566026
$ws.Range("F7").Value = 24.3751090906499
$ws.Range("G7").Value = 23.38682183068499
$ws.Range("H7").Value = 12.94353115243495
$ws.Range("J7").Value = 9.943016321331221
$ws.Range("M7").Value = 53.28508849058028
$ws.Range("O7").Value = 18.98128922860215
$ws.Range("C8").Value = 4.862140570404265
$ws.Range("D8").Value = 9.062114939295229
$ws.Range("E8").Value = 14.77632089438222
$ws.Range("F8").Value = 23.73818810412237
$ws.Range("G8").Value = 22.96594497557528
$ws.Range("H8").Value = 12.76553515916864
$ws.Range("J8").Value = 10.08489798456134
$ws.Range("M8").Value = 57.85209563267954
$ws.Range("O8").Value = 18.66378642252099
$ws.Range("C9").Value = 5.277861113304152
$ws.Range("D9").Value = 8.984543361719473
$ws.Range("E9").Value = 15.47898182267786
$ws.Range("F9").Value = 22.60624768092319
$ws.Range("G9").Value = 22.33996019051571
$ws.Range("H9").Value = 12.45396642972739
$ws.Range("J9").Value = 10.39092549020215
$ws.Range("M9").Value = 65.82832083816206
$ws.Range("O9").Value = 18.12484182855108
$ws.Range("C10").Value = 5.570222006645617
$ws.Range("D10").Value = 8.945521894065966
$ws.Range("E10").Value = 16.00986028328156
$ws.Range("F10").Value = 21.84866941483864
$ws.Range("G10").Value = 22.01447449111392
$ws.Range("H10").Value = 12.24846622997864
$ws.Range("J10").Value = 10.63124776017389
$ws.Range("M10").Value = 71.08854949472322
$ws.Range("O10").Value = 17.78278856228542
$ws.Range("C11").Value = 5.6997439850251
$ws.Range("D11").Value = 8.931829739655749
$ws.Range("E11").Value = 16.25354037498679
$ws.Range("F11").Value = 21.52081394496465
$ws.Range("G11").Value = 21.89839834097976
$ws.Range("H11").Value = 12.16019514680652
$ws.Range("J11").Value = 10.74361500737776
$ws.Range("M11").Value = 73.3515149663338
$ws.Range("O11").Value = 17.63956432107203
$ws.Range("C12").Value = 5.748245450827866
$ws.Range("D12").Value = 8.927241881661603
$ws.Range("E12").Value = 16.3460497673874
$ws.Range("F12").Value = 21.39914499473899
$ws.Range("G12").Value = 21.85927397206713
$ws.Range("H12").Value = 12.1275309442423
$ws.Range("J12").Value = 10.7865756301911
$ws.Range("M12").Value = 74.18980456318577
$ws.Range("O12").Value = 17.58716945365965
$ws.Range("C13").Value = 5.737824737113692
$ws.Range("D13").Value = 8.928203187964147
$ws.Range("E13").Value = 16.32611708482011
$ws.Range("F13").Value = 21.42523692774687
$ws.Range("G13").Value = 21.86748158754227
$ws.Range("H13").Value = 12.1345316577235
$ws.Range("J13").Value = 10.77730551859932
$ws.Range("M13").Value = 74.01009234956204
$ws.Range("O13").Value = 17.59837073570781
$ws.Range("C14").Value = 5.70374541054768
$ws.Range("D14").Value = 8.931440250072392
$ws.Range("E14").Value = 16.2611469535759
$ws.Range("F14").Value = 21.51075396884255
$ws.Range("G14").Value = 21.8950814179094
$ws.Range("H14").Value = 12.15749249331137
$ws.Range("J14").Value = 10.74714135475252
$ws.Range("M14").Value = 73.42085495547067
$ws.Range("O14").Value = 17.63521650451788
$ws.Range("C15").Value = 5.682798500804979
$ws.Range("D15").Value = 8.933501201903237
$ws.Range("E15").Value = 16.22137890527343
$ws.Range("F15").Value = 21.56346122381055
$ws.Range("G15").Value = 21.91262308496066
$ws.Range("H15").Value = 12.17165631129745
$ws.Range("J15").Value = 10.72871744528137
$ws.Range("M15").Value = 73.0575032657287
$ws.Range("O15").Value = 17.65802723163491
$ws.Range("C16").Value = 5.56168345916173
$ws.Range("D16").Value = 8.946499522745818
$ws.Range("E16").Value = 15.99397211163616
$ws.Range("F16").Value = 21.87043889737904
$ws.Range("G16").Value = 22.02272385486442
$ws.Range("H16").Value = 12.25434089511706
$ws.Range("J16").Value = 10.62396313444614
$ws.Range("M16").Value = 70.93804746976106
$ws.Range("O16").Value = 17.79240282520595
$ws.Range("C17").Value = 5.486458836620758
$ws.Range("D17").Value = 8.955522594631093
$ws.Range("E17").Value = 15.85496526488797
$ws.Range("F17").Value = 22.06310307218767
$ws.Range("G17").Value = 22.09862538730496
$ws.Range("H17").Value = 12.30640887630448
$ws.Range("J17").Value = 10.56045923240561
$ws.Range("M17").Value = 69.60456973431124
$ws.Range("O17").Value = 17.87805057260941
$ws.Range("C18").Value = 5.442866975308779
$ws.Range("D18").Value = 8.961093650748261
$ws.Range("E18").Value = 15.77522203582448
$ws.Range("F18").Value = 22.17549218433305
$ws.Range("G18").Value = 22.14527912694761
$ws.Range("H18").Value = 12.33684677275467
$ws.Range("J18").Value = 10.52422160724156
$ws.Range("M18").Value = 68.82535421726817
$ws.Range("O18").Value = 17.92847406348184
$ws.Range("C19").Value = 5.428053257593441
$ws.Range("D19").Value = 8.963044986733379
$ws.Range("E19").Value = 15.74826091430366
$ws.Range("F19").Value = 22.21381347546513
$ws.Range("G19").Value = 22.16158297132772
$ws.Range("H19").Value = 12.34723626760544
$ws.Range("J19").Value = 10.51200251657893
$ws.Range("M19").Value = 68.5594231581964
$ws.Range("O19").Value = 17.94574421268107
$ws.Range("C20").Value = 5.494500577459041
$ws.Range("D20").Value = 8.954522521897072
$ws.Range("E20").Value = 15.86974166193652
$ws.Range("F20").Value = 22.0424302742748
$ws.Range("G20").Value = 22.09023366005661
$ws.Range("H20").Value = 12.30081537230694
$ws.Range("J20").Value = 10.5671897138483
$ws.Range("M20").Value = 69.74778671913535
$ws.Range("O20").Value = 17.86881261410247
$ws.Range("C21").Value = 5.713770484580222
$ws.Range("D21").Value = 8.930473134046498
$ws.Range("E21").Value = 16.28022455117347
$ws.Range("F21").Value = 21.48556754664924
$ws.Range("G21").Value = 21.88684173384202
$ws.Range("H21").Value = 12.15072754993327
$ws.Range("J21").Value = 10.75599040698662
$ws.Range("M21").Value = 73.59443389933207
$ws.Range("O21").Value = 17.62434353973027
$ws.Range("C22").Value = 5.853877771086395
$ws.Range("D22").Value = 8.918242035093453
$ws.Range("E22").Value = 16.54982308396328
$ws.Range("F22").Value = 21.13613684299026
$ws.Range("G22").Value = 21.78217777987578
$ws.Range("H22").Value = 12.05708530570076
$ws.Range("J22").Value = 10.8817572033885
$ws.Range("M22").Value = 75.99978792849865
$ws.Range("O22").Value = 17.47533062995664
$ws.Range("C23").Value = 5.779406322890205
$ws.Range("D23").Value = 8.924446479863349
$ws.Range("E23").Value = 16.40583752280764
$ws.Range("F23").Value = 21.32128091067992
$ws.Range("G23").Value = 21.8353771556796
$ws.Range("H23").Value = 12.10665246222129
$ws.Range("J23").Value = 10.81442499915463
$ws.Range("M23").Value = 74.72592858583477
$ws.Range("O23").Value = 17.55385585156602
$ws.Range("C24").Value = 5.490865975050967
$ws.Range("D24").Value = 8.954973461049908
$ws.Range("E24").Value = 15.8630607057563
$ws.Range("F24").Value = 22.05177139336152
$ws.Range("G24").Value = 22.09401818042031
$ws.Range("H24").Value = 12.30334262868171
$ws.Range("J24").Value = 10.56414601503846
$ws.Range("M24").Value = 69.68307751126132
$ws.Range("O24").Value = 17.87298541352079
$ws.Range("C25").Value = 5.167454195365126
$ws.Range("D25").Value = 9.002424907088912
$ws.Range("E25").Value = 15.28597090918914
$ws.Range("F25").Value = 22.89967831141087
$ws.Range("G25").Value = 22.48664408048192
$ws.Range("H25").Value = 12.53418590203439
$ws.Range("J25").Value = 10.30530873076515
$ws.Range("M25").Value = 63.77585563444067
$ws.Range("O25").Value = 18.26139329634551
